$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13. This shifts existing rows 13..23 down to
# 14..24, carrying their values/styles/heights along automatically (and the
# sheet dimension grows to A1:C24 on its own).
$ws.Rows.Item(13).Insert()

# The freshly inserted row 13 has no B/C cells yet, so grab the number
# format/font/alignment from the (still untouched) B19/C19 pair -- which
# carries exactly the styling every other value cell in columns B/C uses --
# before we overwrite the text content below.
$ws.Range("B19").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C19").Copy()
$ws.Range("C13").PasteSpecial(-4122)

# --- New/changed text content -------------------------------------------------

# Objetivos (row 10): replace placeholder text with the real course objectives.
$objetivos = "1) Gerais: - Levar os estudantes a compreenderem os mecanismos de obtenção da influencia de diversos fatores ( variáveis independentes de um processo) sobre as variáveis resposta ( dependentes), através da análise multivariada.2) Específicos: - Ao final do curso os educandos devem:? Saber planejar e executar um experimento fatorial completo e fracionado? Saber analisar os resultados propondo a condição de melhor ajuste que otimiza os valores da variável resposta na região experimental estudada? Dominar, pelo menos, um software comercial sobre o assunto? Saber modelar um processo , com base em dados empíricos"
$ws.Range("B10").Value = $objetivos
$ws.Range("C10").Value = $objetivos

# Docentes responsáveis (newly inserted row 13): professor name moves here.
$ws.Range("B13").Value = "5840535 - Messias Borges Silva"
$ws.Range("C13").Value = "5840535 - Messias Borges Silva"

# Programa resumido / Programa (rows 14 and 16): new syllabus text.
$programa = "IntroduçãoExperimentação convencionalExperimentos Fatoriais completosExperimentos Fatoriais fracionadosAnálise de variânciaMetodologia de superfície de respostaMétodo de Taguchi"
$ws.Range("B14").Value = $programa
$ws.Range("C14").Value = $programa
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa

# Método (row 19): the course assessment method.
$ws.Range("B19").Value = "2 provas escritas"
$ws.Range("C19").Value = "2 provas escritas"

# Critério (row 20): the assessment criteria.
$criterio = "serão avaliados os conteúdos discutidos em sala e constantes da ementa do curso.A média da disciplina será a média aritmética das duas provas."
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

# Norma de recuperação (row 21): the make-up exam rule.
$ws.Range("B21").Value = "uma provas escrita com conteúdo de todo o semestre"
$ws.Range("C21").Value = "uma provas escrita com conteúdo de todo o semestre"

# Bibliografia (row 22): new bibliography text.
$bibliografia = "1) Planejamento e otimização de Experimentos. Roy E. Bruns, Edit. UNICAMP, 19963) Design and Analysis of Experiments, Douglas C. Montgomery, 6th edition, wiley, 20054) Designing for Quality  Robert H. Lochner  Ed. Quality Press, 19945) Statistics for Experimenter. Box & Hunter"
$ws.Range("B22").Value = $bibliografia
$ws.Range("C22").Value = $bibliografia

# Row 13 has no label in column A (it only carries the professor name in B/C).
$ws.Range("A13").Clear()

# --- Column layout ----------------------------------------------------------
# Column A (labels) and column B (values) now have distinct widths;
# previously columns A:B shared one <col> entry.
$ws.Columns.Item(2).ColumnWidth = 60.7109375
